# Auto-generated Excel COM-interop edit script
# Applies cell value updates to Sheet1 per the target diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.383.65"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.615.63"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.09%  "
$c = $ws.Range("D5")
$c.Value = "'213.29"
$c.Style = "Normal"

$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E6").Value = "  -0.12%  "
$c = $ws.Range("D7")
$c.Value = "'0.486"
$c.Style = "Normal"

$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  -0.70%  "
$c = $ws.Range("D10")
$c.Value = "'18.55"
$c.Style = "Normal"

$ws.Range("E10").Value = "  +2.18%  "
$c = $ws.Range("D11")
$c.Value = "'0.0813"
$c.Style = "Normal"

$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "1.839.39"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "1.629.25"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("D16").Value = "26.382.04"
$ws.Range("E16").Value = "  +0.19%  "
$c = $ws.Range("D17")
$c.Value = "'62.05"
$c.Style = "Normal"

$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  -0.04%  "
$c = $ws.Range("D20")
$c.Value = "'202.96"
$c.Style = "Normal"

$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +5.03%  "
$c = $ws.Range("D25")
$c.Value = "'144.86"
$c.Style = "Normal"

$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("E28").Value = "  -0.41%  "
$c = $ws.Range("D29")
$c.Value = "'6.59"
$c.Style = "Normal"

$ws.Range("E30").Value = "  +4.00%  "
$c = $ws.Range("D31")
$c.Value = "'1.18"
$c.Style = "Normal"

$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  +2.13%  "
$c = $ws.Range("D33")
$c.Value = "'2.95"
$c.Style = "Normal"

$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").Value = "1.164.76"
$ws.Range("E36").Value = "  +4.89%  "
$c = $ws.Range("D37")
$c.Value = "'0.0165"
$c.Style = "Normal"

$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("E38").Value = "  -0.10%  "
$c = $ws.Range("D39")
$c.Value = "'0.794"
$c.Style = "Normal"

$ws.Range("E39").Value = "  +0.76%  "
$c = $ws.Range("D40")
$c.Value = "'2.33"
$c.Style = "Normal"

$ws.Range("E40").Value = "  -0.40%  "
$c = $ws.Range("D41")
$c.Value = "'0.503"
$c.Style = "Normal"

$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("E42").Value = "  +1.07%  "
$c = $ws.Range("D43")
$c.Value = "'5.25"
$c.Style = "Normal"

$ws.Range("D44").Value = "1.753.89"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  -1.51%  "
$c = $ws.Range("D46")
$c.Value = "'1.53"
$c.Style = "Normal"

$ws.Range("E46").Value = "  -1.63%  "
$c = $ws.Range("D47")
$c.Value = "'54.38"
$c.Style = "Normal"

$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D49")
$c.Value = "'0.407"
$c.Style = "Normal"

$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0957"
$ws.Range("E50").Value = "  -9.65%  "
$ws.Range("E51").Value = "  +0.07%  "
